$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 575
$ws.Range("F4").Value = 30
$ws.Range("F5").Value = 738
$ws.Range("F6").Value = 363
$ws.Range("G6").Value = 36
$ws.Range("F10").Value = 215
$ws.Range("F11").Value = 5920
$ws.Range("F12").Value = 51
$ws.Range("F13").Value = 42
$ws.Range("F19").Value = 121
$ws.Range("F21").Value = 705
$ws.Range("F22").Value = 136
$ws.Range("F24").Value = 309
$ws.Range("F25").Value = 1018
$ws.Range("F27").Value = 1803
$ws.Range("F28").Value = 465

# Sheet "演出"
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 300

# Sheet "本地生活"
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 225

# Sheet "全部类型"
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 225
$ws.Range("F3").Value = 576
$ws.Range("F5").Value = 30
$ws.Range("F6").Value = 738
$ws.Range("F8").Value = 363
$ws.Range("G8").Value = 36
$ws.Range("F12").Value = 215
$ws.Range("F13").Value = 5920
$ws.Range("F14").Value = 51
$ws.Range("F15").Value = 42
$ws.Range("F23").Value = 121
$ws.Range("F26").Value = 300
$ws.Range("F28").Value = 705
$ws.Range("F32").Value = 136
$ws.Range("F34").Value = 309
$ws.Range("F35").Value = 1018
$ws.Range("F37").Value = 1803
$ws.Range("F38").Value = 465
